$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.957.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.478.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.26%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.81%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -3.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.479.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.135'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.164'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("E12").Value = '  -4.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.336'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.41%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.12%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.962.12'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.079.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.480.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '362.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("E23").Value = '  -5.13%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("E26").Value = '  -7.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.23%  '
$ws.Range("E28").Value = '  +1.62%  '
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  -7.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.12'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '515.87'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.71%  '
$ws.Range("E33").Value = '  -2.52%  '
$ws.Range("E34").Value = '  -5.88%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  -3.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.35'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  -3.89%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.332'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.29%  '
$ws.Range("E44").Value = '  -3.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("E47").Value = '  -4.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("E49").Value = '  -4.08%  '
$ws.Range("E50").Value = '  -3.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0738'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.71%  '
